$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("B1").Value = 100
$ws.Range("C1").Value = $true

# New row 2 - Lampada da sala
$ws.Range("A2").Value = "Lampada da sala"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = $false

# New row 3 - Ar da sala
$ws.Range("A3").Value = "Ar da sala"
$ws.Range("B3").Value = 18
$ws.Range("C3").Value = $true

# New row 4 - Ar do quarto
$ws.Range("A4").Value = "Ar do quarto"
$ws.Range("B4").Value = 23
$ws.Range("C4").Value = $false
